# Add a new bullet item right after the "Suma odhadů odpovídá
# předpokládané časové kapacitě sprintu" item in the
# "Definition of Ready pro sprint" list.

$d = $word.ActiveDocument

$anchorText = "Suma odhadů odpovídá předpokládané časové kapacitě sprintu"
$newText    = "Pracnost je rovnoměrně rozdělena mezi řešitele v týmu"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r", "`n", "`x07") -eq $anchorText) {
        # Insert a new paragraph right after the anchor; it inherits the
        # anchor's paragraph formatting (ListParagraph style + bullet numPr).
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = $newText
        break
    }
}
